# Apply the "DaySale" shortage-report update:
#  - Insert a new product "COLLOMAK TOP. SOUTION 10 ML" into the sorted
#    product list (it lands at what is currently row 12, pushing every
#    row below it down by one).
#  - Re-number the index ("م") column accordingly.
#  - Recompute the grand total (P column) to include the new line.
#  - Refresh the footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

function Set-ProductRow {
    param($r, $idx, $name, $bal, $ord, $price, $sell, $trans)
    $ws.Cells.Item($r, 1).Value = $idx
    Set-TextCell $ws.Cells.Item($r, 3) $name
    Set-TextCell $ws.Cells.Item($r, 8) $bal
    Set-TextCell $ws.Cells.Item($r, 12) $ord
    Set-TextCell $ws.Cells.Item($r, 14) $price
    Set-TextCell $ws.Cells.Item($r, 16) $sell
    Set-TextCell $ws.Cells.Item($r, 17) $trans
}

# Rows 12-24 shift down by one product each (cell styles stay put - only
# the values need to move along with the new insertion at row 12).
Set-ProductRow 12 6  "COLLOMAK TOP. SOUTION 10 ML"               "1:0"  "1" "28.00"  "28.0000"  "1:0"
Set-ProductRow 13 7  "DANSET 8MG/4ML 3 AMP."                     "0:1"  "1" "142.50" "47.0250"  "0:1"
Set-ProductRow 14 8  "DIASTOP SUSP. 60ML"                        "2:0"  "1" "30.00"  "30.0000"  "1:0"
Set-ProductRow 15 9  "LAMIFEN 1% CREAM 15 GM"                    "2:0"  "1" "18.00"  "18.0000"  "1:0"
Set-ProductRow 16 10 "LANTANON 30MG 10 CAPS."                    "0:0"  "1" "57.00"  "57.0000"  "1:0"
Set-ProductRow 17 11 "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML"     "4:0"  "1" "57.00"  "57.0000"  "1:0"
Set-ProductRow 18 12 "METACARDIA MR 35MG 30 F.C. TAB."           "0:2"  "1" "60.00"  "60.0000"  "1:0"
Set-ProductRow 19 13 "NANAZOXID 500MG 18 F.C. TABS."             "0:1"  "1" "114.00" "37.6200"  "0:1"
Set-ProductRow 20 14 "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML" "3:0"  "1" "24.00"  "24.0000"  "1:0"
Set-ProductRow 21 15 "PROXIMOL COMPOUND EFF. GRANULES 60 GM"     "1:0"  "1" "52.00"  "52.0000"  "1:0"
Set-ProductRow 22 16 "SPASMO-DIGESTIN 30 TABS."                  "2:2"  "1" "78.00"  "78.0000"  "1:0"
Set-ProductRow 23 17 "TELFAST 30MG/5ML SUSP. 100 ML"             "2:0"  "1" "50.00"  "50.0000"  "1:0"
Set-ProductRow 24 18 "حبايه"                                      "0:0"  "0" "3.00"   "3.0000"   "1:0"

# Row 25 used to be the grand-total row; it is recycled into the last
# product line (previously row 24's content), formatted just like the
# other product rows above it.
$ws.Range("A24:Q24").Copy($ws.Range("A25:Q25"))
Set-ProductRow 25 19 "سرنجات 3 سم" "0:0" "0" "2.00" "2.0000" "1:0"

# Row 27 becomes the footer row (was row 26 before the insert); copy the
# formatting down *before* row 26 gets overwritten, then refresh the
# printed timestamp.
$ws.Range("A26:Q26").Copy($ws.Range("A27:Q27"))
$ws.Rows.Item(27).RowHeight = 16.5
$ws.Cells.Item(27, 1).Value = "Thursday, 19 June, 2025 11:55 AM"
$ws.Cells.Item(27, 7).Value = "1/1"
$ws.Cells.Item(27, 11).Value = "developed by : Abdelaziz Talaat"

# Row 26 becomes the new grand-total row (was row 25 before the insert).
# Clear out the footer content it held, copy in the total-row formatting,
# then update the sum and the row height to match the freshly generated
# report.
$ws.Range("A26:Q26").ClearContents()
$ws.Range("P25:Q25").Copy($ws.Range("P26:Q26"))
$ws.Cells.Item(26, 16).Value = 1203.375
$ws.Rows.Item(26).RowHeight = 25.5

# Fix up the merged-cell regions: the old P25:Q25 / A26:F26 / G26:I26 /
# K26:Q26 merges need to move down to match the new row layout.
$ws.Range("P25:Q25").UnMerge()
$ws.Range("A26:F26").UnMerge()
$ws.Range("G26:I26").UnMerge()
$ws.Range("K26:Q26").UnMerge()

$ws.Range("A25:B25").Merge()
$ws.Range("C25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()
$ws.Range("N25:O25").Merge()

$ws.Range("P26:Q26").Merge()
$ws.Range("A27:F27").Merge()
$ws.Range("G27:I27").Merge()
$ws.Range("K27:Q27").Merge()
